$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.394.67"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.47"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.31"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6295"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07623"
$ws.Range("E8").Value = "  +1.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2940"
$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.47"
$ws.Range("E10").Value = "  +0.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07747"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.843.71"
$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.008"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001088"
$ws.Range("E14").Value = "  +8.76%  "

$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.44"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.091.44"
$ws.Range("E17").Value = "  -7.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.134"
$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.426.67"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.41"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.440"
$ws.Range("E23").Value = "  -1.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.30"
$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("E26").Value = "  -0.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.375"
$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.65"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.466"
$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("E30").Value = "  +3.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05629"
$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.112"
$ws.Range("E32").Value = "  -0.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.041"
$ws.Range("E33").Value = "  +0.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.852"
$ws.Range("E34").Value = "  +0.44%  "

$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7094"
$ws.Range("E36").Value = "  -0.72%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.776"
$ws.Range("E38").Value = "  -0.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.229.86"
$ws.Range("E39").Value = "  -1.95%  "

$ws.Range("E40").Value = "  -0.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.486"
$ws.Range("E41").Value = "  +4.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9078"
$ws.Range("E42").Value = "  -0.59%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.45"
$ws.Range("E44").Value = "  +0.40%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.04"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000122"
$ws.Range("E46").Value = "  +3.05%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.151"
$ws.Range("E47").Value = "  +1.49%  "

$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4010"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.026"
$ws.Range("E49").Value = "  -0.96%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.685"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1121"
$ws.Range("E51").Value = "  -0.57%  "
